$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '57.866.60'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.92%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.452.81'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.20%  '
$ws.Range("E4").Value = '  +0.23%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '526.50'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.01%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '130.64'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.23%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.98%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.563'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.21%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.450.84'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.68%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0979'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.33%  '
$ws.Range("E11").Value = '  -3.22%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.94'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.18%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.321'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.97%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.892.38'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.96%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '57.779.81'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.83%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.76'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.70%  '
$ws.Range("E17").Value = '  -1.60%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.463.78'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.07%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.34'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.66%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.15'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.36%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '316.53'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.49%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.06'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.47%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.34%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.94'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.70%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.406'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.30%  '
$ws.Range("E26").Value = '  +1.09%  '
$ws.Range("E27").Value = '  -2.50%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.28'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.02%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '171.79'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.59%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0733'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.13%  '
$ws.Range("E31").Value = '  -0.56%  '
$ws.Range("E32").Value = '  -3.83%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.09'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.66%  '
$ws.Range("E34").Value = '  +0.21%  '
$ws.Range("E35").Value = '  +0.55%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.78'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.54%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.18'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -6.56%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.80'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.00%  '
$ws.Range("B39").Value = 'OKB'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.20'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.66%  '
$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.46'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.01%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.800'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.95%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.40'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.02%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '267.19'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.72%  '
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.84'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.85%  '
$ws.Range("B45").Value = 'Mantle'
$ws.Range("C45").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.583'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.58%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '125.26'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.62%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0930'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.12%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0494'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.10%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0211'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.70%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.01'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.37%  '
$ws.Range("B51").Value = 'InjectiveProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '16.34'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.47%  '
